$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 289401.03
$ws.Range("J17").Value = 306839.7
$ws.Range("L17").Value = 920519.1000000001
$ws.Range("N17").Value = -920855.1000000001

$ws.Range("H53").Value = 689.55
$ws.Range("I53").Value = 145.5
$ws.Range("J53").Value = 825.5625
$ws.Range("K53").Value = 145.5
$ws.Range("L53").Value = 825.5625
$ws.Range("M53").Value = 491.5
$ws.Range("N53").Value = -2099.5625

$ws.Range("H70").Value = 8730.074000000001
$ws.Range("J70").Value = 10200.632
$ws.Range("L70").Value = 30601.896
$ws.Range("N70").Value = -31141.896

$ws.Range("H73").Value = 8730.074000000001
$ws.Range("J73").Value = 10200.632
$ws.Range("L73").Value = 30601.896
$ws.Range("N73").Value = -32473.896

$ws.Range("H106").Value = 3396.9443
$ws.Range("I106").Value = 2390.6365
$ws.Range("K106").Value = 2390.6365
$ws.Range("M106").Value = -1759.6365

$ws.Range("H111").Value = 26222.455
$ws.Range("J111").Value = 3796.6667
$ws.Range("L111").Value = 11390.0001
$ws.Range("N111").Value = -17524.0001

$ws.Range("H112").Value = 272868.22
$ws.Range("J112").Value = 288410.7
$ws.Range("L112").Value = 865232.1000000001
$ws.Range("N112").Value = -867448.1000000001

$ws.Range("H125").Value = 3342.6
$ws.Range("I125").Value = 1170.8
$ws.Range("K125").Value = 10537.2
$ws.Range("M125").Value = -8077.199999999999

$ws.Range("H129").Value = 3598
$ws.Range("I129").Value = 643.2308
$ws.Range("K129").Value = 1929.6924
$ws.Range("M129").Value = 3070.3076

$ws.Range("H132").Value = 52636700
$ws.Range("I132").Value = 55560684
$ws.Range("K132").Value = 166682052
$ws.Range("M132").Value = -166679522

$ws.Range("H133").Value = 44875
$ws.Range("J133").Value = 44875
$ws.Range("L133").Value = 44875
$ws.Range("N133").Value = -54995

$ws.Range("H137").Value = 2511.8286
$ws.Range("I137").Value = 2586.9565
$ws.Range("J137").Value = 2367.8333
$ws.Range("K137").Value = 7760.869499999999
$ws.Range("L137").Value = 7103.499899999999
$ws.Range("M137").Value = -5210.869499999999
$ws.Range("N137").Value = -12203.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 12000
$ws.Range("J11").Value = 12000
$ws.Range("L11").Value = 12000
$ws.Range("N11").Value = -12288

$ws.Range("H32").Value = 1612.6224
$ws.Range("I32").Value = 1612.6224
$ws.Range("K32").Value = 1612.6224
$ws.Range("M32").Value = -1325.6224

$ws.Range("H61").Value = 854915.75
$ws.Range("I61").Value = 1255338.5
$ws.Range("J61").Value = 11920.579
$ws.Range("K61").Value = 1255338.5
$ws.Range("L61").Value = 11920.579
$ws.Range("M61").Value = -1255126.5
$ws.Range("N61").Value = -12344.579

$ws.Range("H102").Value = 6649.9
$ws.Range("J102").Value = 4999
$ws.Range("L102").Value = 4999
$ws.Range("N102").Value = -8243

$ws.Range("H110").Value = 1906.25
$ws.Range("I110").Value = 1875
$ws.Range("K110").Value = 1875
$ws.Range("M110").Value = 170

$ws.Range("H122").Value = 4310.4
$ws.Range("I122").Value = 4068.138
$ws.Range("J122").Value = 5481.3335
$ws.Range("K122").Value = 12204.414
$ws.Range("L122").Value = 16444.0005
$ws.Range("M122").Value = -9754.414000000001
$ws.Range("N122").Value = -21344.0005

$ws.Range("H136").Value = 854915.75
$ws.Range("I136").Value = 1255338.5
$ws.Range("J136").Value = 11920.579
$ws.Range("K136").Value = 3766015.5
$ws.Range("L136").Value = 35761.737
$ws.Range("M136").Value = -3763465.5
$ws.Range("N136").Value = -40861.737

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2009.175
$ws.Range("I20").Value = 1741.303
$ws.Range("J20").Value = 3272
$ws.Range("K20").Value = 1741.303
$ws.Range("L20").Value = 3272
$ws.Range("M20").Value = -1494.303
$ws.Range("N20").Value = -3766

$ws.Range("H99").Value = 4060.6667
$ws.Range("I99").Value = 3449.5715
$ws.Range("K99").Value = 3449.5715
$ws.Range("M99").Value = -1951.5715

$ws.Range("H107").Value = 200
$ws.Range("I107").Value = 200
$ws.Range("K107").Value = 200
$ws.Range("M107").Value = 1720

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2365.8572
$ws.Range("I31").Value = 2365.8572
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2365.8572
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -2070.8572

$ws.Range("H34").Value = 2365.8572
$ws.Range("I34").Value = 2365.8572
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2365.8572
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -2163.8572

$ws.Range("H58").Value = 2702.8518
$ws.Range("J58").Value = 3334.2
$ws.Range("L58").Value = 3334.2
$ws.Range("N58").Value = -3740.2

$ws.Range("H86").Value = 5819.909
$ws.Range("J86").Value = 7010.1665
$ws.Range("L86").Value = 7010.1665
$ws.Range("N86").Value = -9256.166499999999

$ws.Range("H89").Value = 5819.909
$ws.Range("J89").Value = 7010.1665
$ws.Range("L89").Value = 35050.8325
$ws.Range("N89").Value = -46282.8325

$ws.Range("H107").Value = 823.34784
$ws.Range("I107").Value = 476.35715
$ws.Range("K107").Value = 476.35715
$ws.Range("M107").Value = 1443.64285

$ws.Range("H132").Value = 933022.0600000001
$ws.Range("I132").Value = 1741255.1
$ws.Range("K132").Value = 5223765.300000001
$ws.Range("M132").Value = -5221235.300000001

$ws.Range("H134").Value = 4898.4053
$ws.Range("I134").Value = 1220.0588
$ws.Range("K134").Value = 3660.1764
$ws.Range("M134").Value = -1125.1764

$ws.Range("H136").Value = 2702.8518
$ws.Range("J136").Value = 3334.2
$ws.Range("L136").Value = 10002.6
$ws.Range("N136").Value = -15102.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 132.30302
$ws.Range("I2").Value = 142.04762
$ws.Range("J2").Value = 115.25
$ws.Range("K2").Value = 852.28572
$ws.Range("L2").Value = 691.5
$ws.Range("M2").Value = -739.28572
$ws.Range("N2").Value = -917.5

$ws.Range("H11").Value = 352.5
$ws.Range("J11").Value = 300
$ws.Range("L11").Value = 900
$ws.Range("N11").Value = -1180

$ws.Range("H12").Value = 1930.8182
$ws.Range("I12").Value = 3101.2222
$ws.Range("J12").Value = 1120.5385
$ws.Range("K12").Value = 9303.6666
$ws.Range("L12").Value = 3361.6155
$ws.Range("M12").Value = -9130.6666
$ws.Range("N12").Value = -3707.6155

$ws.Range("H55").Value = 10518.091
$ws.Range("J55").Value = 11499.947
$ws.Range("L55").Value = 34499.841
$ws.Range("N55").Value = -34853.841

$ws.Range("H118").Value = 3643
$ws.Range("I118").Value = 464.5
$ws.Range("K118").Value = 1393.5
$ws.Range("M118").Value = -150.5

$ws.Range("H131").Value = 228769.05
$ws.Range("J131").Value = 1707.8334
$ws.Range("L131").Value = 5123.5002
$ws.Range("N131").Value = -15203.5002

$ws.Range("H134").Value = 3354.5
$ws.Range("I134").Value = 2264.3333
$ws.Range("J134").Value = 6625
$ws.Range("K134").Value = 6792.999899999999
$ws.Range("L134").Value = 19875
$ws.Range("M134").Value = -1722.999899999999
$ws.Range("N134").Value = -30015

$ws.Range("H140").Value = 4238.8184
$ws.Range("I140").Value = 3754.6155
$ws.Range("J140").Value = 4938.222
$ws.Range("K140").Value = 11263.8465
$ws.Range("L140").Value = 14814.666
$ws.Range("M140").Value = -6083.8465
$ws.Range("N140").Value = -25174.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5101.357
$ws.Range("I70").Value = 4402.6665
$ws.Range("J70").Value = 5625.375
$ws.Range("K70").Value = 4402.6665
$ws.Range("L70").Value = 5625.375
$ws.Range("M70").Value = -4132.6665
$ws.Range("N70").Value = -6165.375

$ws.Range("H73").Value = 5101.357
$ws.Range("I73").Value = 4402.6665
$ws.Range("J73").Value = 5625.375
$ws.Range("K73").Value = 4402.6665
$ws.Range("L73").Value = 5625.375
$ws.Range("M73").Value = -3466.6665
$ws.Range("N73").Value = -7497.375

$ws.Range("H97").Value = 1922.48
$ws.Range("I97").Value = 985.9286
$ws.Range("K97").Value = 985.9286
$ws.Range("M97").Value = -489.9286

$ws.Range("H102").Value = 19619.484
$ws.Range("I102").Value = 3391
$ws.Range("K102").Value = 3391
$ws.Range("M102").Value = -1769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 6540.2383
$ws.Range("I82").Value = 11180
$ws.Range("J82").Value = 2322.2727
$ws.Range("K82").Value = 11180
$ws.Range("L82").Value = 2322.2727
$ws.Range("M82").Value = -10819
$ws.Range("N82").Value = -3044.2727

$ws.Range("H85").Value = 6540.2383
$ws.Range("I85").Value = 11180
$ws.Range("J85").Value = 2322.2727
$ws.Range("K85").Value = 11180
$ws.Range("L85").Value = 2322.2727
$ws.Range("M85").Value = -9932
$ws.Range("N85").Value = -4818.2727

$ws.Range("H122").Value = 7454.5
$ws.Range("I122").Value = 6957.143
$ws.Range("J122").Value = 8150.8
$ws.Range("K122").Value = 20871.429
$ws.Range("L122").Value = 24452.4
$ws.Range("M122").Value = -18421.429
$ws.Range("N122").Value = -29352.4

$ws.Range("H125").Value = 104998.8
$ws.Range("J125").Value = 104998.8
$ws.Range("L125").Value = 104998.8
$ws.Range("N125").Value = -114838.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 21000
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

$ws.Range("H126").Value = 2620.5386
$ws.Range("I126").Value = 2409.5
$ws.Range("J126").Value = 2958.2
$ws.Range("K126").Value = 7228.5
$ws.Range("L126").Value = 8874.599999999999
$ws.Range("M126").Value = -4758.5
$ws.Range("N126").Value = -13814.6
